# Daily update for the water-delivery tracking sheet.
#
# Column layout: A=行号 B=店铺名称 C=地址 D=总天(total days) E=剩余(days remaining)
#                F=开始时间(start date, yyyyMMdd) G/H/I=备注(notes)
#
# Each day that passes, the "days remaining" (E) counter for every shop
# ticks down by one. When a shop's counter would hit zero (i.e. it was at 1
# remaining), the water has been refilled/restocked instead: the counter
# resets to a full cycle of 10 and the start date (F) advances by 10 days.
#
# One row has a corrupted (non 8-digit) start date and is left untouched,
# matching the source data's behavior of skipping rows whose date can't be
# parsed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null) { continue }
    if ($fVal -eq $null) { continue }

    # Skip rows whose start date isn't a clean 8-digit yyyyMMdd value.
    $fStr = [string]$fVal
    if ($fStr.Length -ne 8) { continue }

    if ($eVal -eq 1) {
        # Out of water today -> refill: reset remaining days and roll the
        # start date forward by a full 10-day cycle.
        $eCell.Value = 10
        $fCell.Value = $fVal + 10
    } else {
        # Otherwise just count down one more day.
        $eCell.Value = $eVal - 1
    }
}
